$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": row 2 (columns B..L) ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.5
$summary.Range("C2").Value = 0.5
$summary.Range("D2").Value = 1
$summary.Range("E2").Value = 0.6666666666666666
$summary.Range("F2").Value = 0.8333333333333334
$summary.Range("G2").Value = 0.9629629629629629
$summary.Range("H2").Value = 0.7774832021770539
$summary.Range("I2").Value = 534
$summary.Range("J2").Value = 534
$summary.Range("K2").Value = 0
$summary.Range("L2").Value = 0

# --- Sheet "Classification Report" ---
$report = $wb.Worksheets.Item("Classification Report")

# Row 2: class "0"
$report.Range("B2").Value = 0
$report.Range("C2").Value = 0
$report.Range("D2").Value = 0

# Row 3: class "1"
$report.Range("B3").Value = 0.5
$report.Range("C3").Value = 1
$report.Range("D3").Value = 0.6666666666666666

# Row 4: accuracy
$report.Range("B4").Value = 0.5
$report.Range("C4").Value = 0.5
$report.Range("D4").Value = 0.5
$report.Range("E4").Value = 0.5

# Row 5: macro avg
$report.Range("B5").Value = 0.25
$report.Range("C5").Value = 0.5
$report.Range("D5").Value = 0.3333333333333333

# Row 6: weighted avg
$report.Range("B6").Value = 0.25
$report.Range("C6").Value = 0.5
$report.Range("D6").Value = 0.3333333333333333

# --- Sheet "Confusion Matrix" ---
$cm = $wb.Worksheets.Item("Confusion Matrix")
$cm.Range("B2").Value = 0
$cm.Range("C2").Value = 534
$cm.Range("B3").Value = 0
$cm.Range("C3").Value = 534
